$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 3,4,5,8,9,10,11,12
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 7
$ws.Range("F5").Value = -5
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = -7
$ws.Range("F10").Value = 6
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = -3
